# Update countries & provincias Spain
#
# Source data is ranked by "Casos totales" (col B) descending. On this
# update: Francia overtakes Iran, Peru overtakes Eslovenia/Argentina/
# Croacia/Mexico/Republica Dominicana, and Afganistan overtakes Venezuela/
# Sri Lanka/Camboya/Costa de Marfil/Mauricio, so those rows re-sort and
# pick up fresh case numbers for the country that moved up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4) : updated totals ---------------------------
$ws.Range("B4").Value = 96968
$ws.Range("C4").Value = 11533
$ws.Range("E4").Value = 93038
$ws.Range("G4").Value = 182
$ws.Range("H4").Value = 1477

# --- Alemania (row 8) : updated totals ----------------------------------
$ws.Range("B8").Value = 50178
$ws.Range("C8").Value = 6240
$ws.Range("E8").Value = 44184

# --- Francia overtakes Iran (rows 9-10) ---------------------------------
$ws.Range("A9").Value = "Francia"
$ws.Range("B9").Value = 32964
$ws.Range("C9").Value = 3809
$ws.Range("D9").Value = 5700
$ws.Range("E9").Value = 25269
$ws.Range("F9").Value = 3787
$ws.Range("G9").Value = 299
$ws.Range("H9").Value = 1995

$ws.Range("A10").Value = "Iran"
$ws.Range("B10").Value = 32332
$ws.Range("C10").Value = 2926
$ws.Range("D10").Value = 11133
$ws.Range("E10").Value = 18821
$ws.Range("F10").Value = 2893
$ws.Range("G10").Value = 144
$ws.Range("H10").Value = 2378

# --- Austria (row 15) : updated totals ----------------------------------
$ws.Range("B15").Value = 7642
$ws.Range("C15").Value = 733
$ws.Range("E15").Value = 7359

# --- Rumania (row 35) : updated totals ----------------------------------
$ws.Range("E35").Value = 1151
$ws.Range("G35").Value = 3
$ws.Range("H35").Value = 26

# --- Peru overtakes Eslovenia, Argentina, Croacia, Mexico and
#     Republica Dominicana (rows 49-54) ----------------------------------
$ws.Range("A49").Value = "Peru"
$ws.Range("B49").Value = 635
$ws.Range("C49").Value = 55
$ws.Range("D49").Value = 14
$ws.Range("E49").Value = 612
$ws.Range("F49").Value = 14
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 9

$ws.Range("A50").Value = "Eslovenia"
$ws.Range("B50").Value = 632
$ws.Range("C50").Value = 70
$ws.Range("D50").Value = 10
$ws.Range("E50").Value = 613
$ws.Range("F50").Value = 14
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 9

$ws.Range("A51").Value = "Argentina"
$ws.Range("B51").Value = 589
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 72
$ws.Range("E51").Value = 504
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 13

$ws.Range("A52").Value = "Croacia"
$ws.Range("B52").Value = 586
$ws.Range("C52").Value = 91
$ws.Range("D52").Value = 37
$ws.Range("E52").Value = 546
$ws.Range("F52").Value = 14
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 3

$ws.Range("A53").Value = "Mexico"
$ws.Range("B53").Value = 585
$ws.Range("C53").Value = 110
$ws.Range("D53").Value = 4
$ws.Range("E53").Value = 573
$ws.Range("F53").Value = 1
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 8

$ws.Range("A54").Value = "Republica Dominicana"
$ws.Range("B54").Value = 581
$ws.Range("C54").Value = 93
$ws.Range("D54").Value = 3
$ws.Range("E54").Value = 558
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 10
$ws.Range("H54").Value = 20

# --- Afganistan overtakes Venezuela, Sri Lanka, Camboya, Costa de
#     Marfil and Mauricio (rows 99-104) -----------------------------------
$ws.Range("A99").Value = "Afganistan"
$ws.Range("B99").Value = 110
$ws.Range("C99").Value = 16
$ws.Range("D99").Value = 2
$ws.Range("E99").Value = 104
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 4

$ws.Range("A100").Value = "Venezuela"
$ws.Range("B100").Value = 107
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 31
$ws.Range("E100").Value = 75
$ws.Range("F100").Value = 2
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 1

$ws.Range("A101").Value = "Sri Lanka"
$ws.Range("B101").Value = 106
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 7
$ws.Range("E101").Value = 99
$ws.Range("F101").Value = 5
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0

$ws.Range("A102").Value = "Camboya"
$ws.Range("B102").Value = 99
$ws.Range("C102").Value = 1
$ws.Range("D102").Value = 11
$ws.Range("E102").Value = 88
$ws.Range("F102").Value = 1
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 0

$ws.Range("A103").Value = "Costa de Marfil"
$ws.Range("B103").Value = 96
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 3
$ws.Range("E103").Value = 93
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 0

$ws.Range("A104").Value = "Mauricio"
$ws.Range("B104").Value = 94
$ws.Range("C104").Value = 13
$ws.Range("D104").Value = 0
$ws.Range("E104").Value = 92
$ws.Range("F104").Value = 1
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 2
